$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)
    $cell.Value = "'2008-04-27"
    $cell.Style = "Normal"
}
